# 09/02 - TA - 3601 DS - Done. Misses PC.
# Adds a new "D03 Variables" sheet (mirroring the layout/styling of the
# existing "D03NonQuotableProducts" sheet) holding the
# textExistingBillingAccountIdField variable, and updates the stored
# selection on a couple of other sheets.

$wb = $excel.ActiveWorkbook

# --- Update selection on "Environments_OnGoing": B22 -> A23 -------------
$wsEnv = $wb.Worksheets.Item("Environments_OnGoing")
$wsEnv.Activate() | Out-Null
$wsEnv.Range("A23").Select() | Out-Null

# --- Update selection on "VoiceContinuity": B3 -> A2 ---------------------
$wsVC = $wb.Worksheets.Item("VoiceContinuity")
$wsVC.Activate() | Out-Null
$wsVC.Range("A2").Select() | Out-Null

# --- Add the new "D03 Variables" sheet at the very end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "D03 Variables"

# Header row, matching the other "Variables" sheets (bold, "Variable"/"Value")
$newSheet.Range("A1").Value = "Variable"
$newSheet.Range("B1").Value = "Value"
$newSheet.Range("A1:B1").Font.Bold = $true

# Data row
$newSheet.Range("A2").Value = "textExistingBillingAccountIdField"
$newSheet.Range("B2").Value = 4121986

# A couple of blank, text-formatted placeholder rows beneath (mirrors the
# other variable sheets in this workbook)
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").HorizontalAlignment = -4131
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").HorizontalAlignment = -4131

# Column widths matching the sibling sheets
$newSheet.Columns.Item(1).ColumnWidth = 31.333333333333332
$newSheet.Columns.Item(2).ColumnWidth = 89.16666666666667

# Page setup matching the sibling sheets
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Leave the new sheet active with A2 selected, and as the selected tab
$newSheet.Activate() | Out-Null
$newSheet.Range("A2").Select() | Out-Null
